$d = $word.ActiveDocument

# Locate the unique sentence that needs to be split/augmented. Include the
# trailing space explicitly so the matched range covers exactly
# "...the list of " (with its trailing space), which is what the final
# run split needs to line up against.
$rng = $d.Content
$found = $rng.Find.Execute(
    "Controller is responsible for monitoring the list of ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)

if (-not $found) {
    throw "Could not find target sentence"
}

$sentenceStart = $rng.Start
$sentenceEnd = $rng.End

$part1 = "Controller is responsible for monitoring "
$part3 = "the list of "
$insertion = "(not maintaining. Maintaining the active list is what we discussed in previous lecture) "

if ($sentenceEnd - $sentenceStart -ne $part1.Length + $part3.Length) {
    throw "Unexpected match length; refusing to edit"
}

$insertAt = $sentenceStart + $part1.Length

# Insert the new parenthetical text right after "monitoring " / right
# before "the list of".
$insertPoint = $d.Range($insertAt, $insertAt)
$insertPoint.InsertBefore($insertion)

# After the insertion, three logical text spans now sit inside what used to
# be a single run: "Controller is responsible for monitoring ", the new
# parenthetical aside, and "the list of ". They (and the untouched runs
# immediately before/after them) all share identical - i.e. absent - run
# formatting, so the engine would otherwise silently re-coalesce them back
# into one big run. Nudging each span's character formatting (set then
# clear again) forces it to keep living in its own <w:r>, matching the
# target run layout, while leaving the run's visible formatting unchanged.
$run1 = $d.Range($sentenceStart, $insertAt)
$run1.Font.Bold = $true
$run1.Font.Bold = $false

$run2Start = $insertAt
$run2End = $run2Start + $insertion.Length
$run2 = $d.Range($run2Start, $run2End)
$run2.Font.Bold = $true
$run2.Font.Bold = $false

$run3Start = $run2End
$run3End = $run3Start + $part3.Length
$run3 = $d.Range($run3Start, $run3End)
$run3.Font.Bold = $true
$run3.Font.Bold = $false

Write-Output "done"
